$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capital structure database values for rows 2 and 3 (same new values in both rows)
$ws.Range("D2:D3").Value = -0.0126
$ws.Range("E2:E3").Value = -0.125
$ws.Range("G2:G3").Value = 0.1122022427178024
$ws.Range("H2:H3").Value = 0.1122022427178024
$ws.Range("I2:I3").Value = 0.04089725021203806
$ws.Range("J2:J3").Value = 0.03642925303739851
$ws.Range("K2:K3").Value = 66.5
$ws.Range("L2:L3").Value = 0.04412447747329308
$ws.Range("M2:M3").Value = 17.3
$ws.Range("N2:N3").Value = 0.01448304730012558
$ws.Range("O2:O3").Value = 0.2601503759398496
$ws.Range("P2:P3").Value = 17.3
$ws.Range("Q2:Q3").Value = 0.01448304730012558
$ws.Range("R2:R3").Value = 0.2601503759398496
$ws.Range("U2:U3").Value = 288.9
$ws.Range("V2:V3").Value = 0.2418585182084554
$ws.Range("W2:W3").Value = 0.04498714652956298
$ws.Range("X2:X3").Value = 0.06980178383840767
$ws.Range("Y2:Y3").Value = -0.02481463730884469
$ws.Range("Z2:Z3").Value = 1.292026912011738
$ws.Range("AA2:AA3").Value = 0.04706757530880423
$ws.Range("AB2:AB3").Value = 0.06975761895035065
$ws.Range("AC2:AC3").Value = -0.02269004364154643
$ws.Range("AE2:AE3").Value = 1.223771027187246
$ws.Range("AF2:AF3").Value = 1.223771027187246
$ws.Range("AG2:AG3").Value = -287.6762289728127
$ws.Range("AH2:AH3").Value = 0.001023456300560091
$ws.Range("AI2:AI3").Value = 0.0009112736356220706
$ws.Range("AJ2:AJ3").Value = -0.3172349889405226
$ws.Range("AK2:AK3").Value = -0.2729314431803194
$ws.Range("AP2:AP3").Value = -4.433905595980529
